# Updated Diversion to account for possible no jail.
# Appends 15 new case rows (253-267) to Sheet1, mirroring the pattern of
# the existing data rows (columns A-K, all plain text values).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$rows = @(
    @("21CRB01291","Bunner","PERMISSION REQ'D TO USE LICENSED DOCK","1501:46-12-04","MM","No Contest","Guilty","$ 500","$ 0","None","None"),
    @("21CRB01291","Bunner","PERMISSION REQ'D TO USE LICENSED DOCK","1501:46-12-04","MM","No Contest","Guilty","$ 34","$ 0","None","None"),
    @("21CRB01291","Bunner","PERMISSION REQ'D TO USE LICENSED DOCK","1501:46-12-04","MM","No Contest","Guilty","$ 34","$ 0","None","None"),
    @("21TRD09437","Bunner","DUS","4510.11","M1","Guilty","Guilty","$ 0","$ 0","None","None"),
    @("21TRD09437","Bunner","1ST SPEED 1 YR SCHOOL >35MPHM4","4511.21B1A","M4","Guilty","Guilty","$ 0","$ 0","None","None"),
    @("21TRD09437","Bunner","RECKLESS OPERATION 1ST IN 1 YR","4511.20","MM","Guilty","Guilty","$ 0","$ 0","None","None"),
    @("21TRD09437","Bunner","DUS","4510.11","M1","Guilty","Guilty","$ 0","$ 0","None","None"),
    @("21TRD09437","Bunner","1ST SPEED 1 YR SCHOOL >35MPHM4","4511.21B1A","M4","Guilty","Guilty","$ 0","$ 0","None","None"),
    @("21TRD09437","Bunner","RECKLESS OPERATION 1ST IN 1 YR","4511.20","MM","Guilty","Guilty","$ 0","$ 0","None","None"),
    @("21CRB01437","Bunner","POSSESSION OF MARIHUANA","2925.11C3","MM","No Contest","Guilty","$ 0","$ 0","None","None"),
    @("03TRD13906","Bunner","SPEED REDUCED ZONE 1ST OFFENSE","4511.21C*","MM","No Contest","Guilty","$ 0","$ 0","None","None"),
    @("21CRB01437","Bunner","POSSESSION OF MARIHUANA","2925.11C3","MM","Guilty","Guilty","$ 0","$ 0","None","None"),
    @("21CRB01437","Bunner","POSSESSION OF MARIHUANA","2925.11C3","MM","Guilty","Guilty","$ 0","$ 0","None","None"),
    @("21CRB01437","Bunner","POSSESSION OF MARIHUANA","2925.11C3","MM","No Contest","Guilty","$ 0","$ 0","None","None"),
    @("21CRB01291","Hemmeter","PERMISSION REQ'D TO USE LICENSED DOCK","1501:46-12-04","MM","Guilty","Guilty","$ 0","$ 0","None","None")
)

$startRow = 253
for ($i = 0; $i -lt $rows.Count; $i++) {
    $r = $startRow + $i
    $values = $rows[$i]
    for ($c = 1; $c -le $values.Count; $c++) {
        $cell = $ws.Cells.Item($r, $c)
        # Force every value to be stored as literal text (matches the
        # existing sheet convention where "$ 0"-style amounts and case
        # numbers like "4511.21C*" are plain text, not numbers/currency).
        $cell.Value = "'" + $values[$c - 1]
        $cell.ClearFormats()
    }
}
